# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row (row 1) - new columns AC, AD, AE
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy the header style used by the existing header row (e.g. A1) onto the
# newly added header cells so they match the rest of the header formatting.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AC1:AE1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill in the team record for every data row (rows 2 through 43).
$lastRow = 43
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 75   # AC - Wins
    $ws.Cells.Item($r, 30).Value = 87   # AD - Losses
    $ws.Cells.Item($r, 31).Value = 0    # AE - Ties
}
